$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 4).Value = '26.779.00'
$ws.Cells.Item(2, 5).Value = '  +4.43%  '

$ws.Cells.Item(3, 4).Value = '1.865.98'
$ws.Cells.Item(3, 5).Value = '  +2.91%  '

$c = $ws.Cells.Item(4, 4)
$c.Value = "'1.001"
$c.Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  +0.24%  '

$c = $ws.Cells.Item(5, 4)
$c.Value = "'269.22"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -2.69%  '

$c = $ws.Cells.Item(6, 4)
$c.Value = "'1.000"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +0.16%  '

$c = $ws.Cells.Item(7, 4)
$c.Value = "'0.5291"
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  +4.77%  '

$c = $ws.Cells.Item(8, 4)
$c.Value = "'0.3355"
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  -2.27%  '

$ws.Cells.Item(9, 2).Value = 'Dogecoin'
$ws.Cells.Item(9, 3).Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$c = $ws.Cells.Item(9, 4)
$c.Value = "'0.06787"
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  +0.78%  '

$ws.Cells.Item(10, 2).Value = 'Solana'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$c = $ws.Cells.Item(10, 4)
$c.Value = "'19.67"
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +0.18%  '

$ws.Cells.Item(11, 2).Value = 'Polygon'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$c = $ws.Cells.Item(11, 4)
$c.Value = "'0.7852"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -2.26%  '

$ws.Cells.Item(12, 2).Value = 'TRON'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$c = $ws.Cells.Item(12, 4)
$c.Value = "'0.07752"
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  -1.04%  '

$ws.Cells.Item(13, 2).Value = 'WrappedEther'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(13, 4).Value = '1.857.08'
$ws.Cells.Item(13, 5).Value = '  +2.46%  '

$ws.Cells.Item(14, 2).Value = 'Litecoin'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$c = $ws.Cells.Item(14, 4)
$c.Value = "'89.71"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +2.59%  '

$ws.Cells.Item(15, 2).Value = 'Polkadot'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Cells.Item(15, 4)
$c.Value = "'5.107"
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +1.24%  '

$ws.Cells.Item(16, 2).Value = 'BinanceUSD'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$c = $ws.Cells.Item(16, 4)
$c.Value = "'1.001"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +0.43%  '

$ws.Cells.Item(17, 2).Value = 'Avalanche'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$c = $ws.Cells.Item(17, 4)
$c.Value = "'14.38"
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  +2.45%  '

$ws.Cells.Item(18, 2).Value = 'ShibaInu'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$c = $ws.Cells.Item(18, 4)
$c.Value = "'0.000008003"
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  +0.08%  '

$ws.Cells.Item(19, 2).Value = 'Dai'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Cells.Item(19, 4)
$c.Value = "'1.001"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +0.20%  '

$ws.Cells.Item(20, 2).Value = 'WrappedBTC'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(20, 4).Value = '26.790.86'
$ws.Cells.Item(20, 5).Value = '  +4.27%  '

$ws.Cells.Item(21, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(21, 4).Value = '2.090.14'
$ws.Cells.Item(21, 5).Value = '  +2.22%  '

$ws.Cells.Item(22, 2).Value = 'Uniswap'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$c = $ws.Cells.Item(22, 4)
$c.Value = "'4.655"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -1.41%  '

$ws.Cells.Item(23, 2).Value = 'Cosmos'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Cells.Item(23, 4)
$c.Value = "'9.880"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -0.45%  '

$ws.Cells.Item(24, 2).Value = 'Chainlink'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Cells.Item(24, 4)
$c.Value = "'6.054"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  -1.16%  '

$ws.Cells.Item(25, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Cells.Item(25, 4)
$c.Value = "'2.394"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +4.72%  '

$ws.Cells.Item(26, 2).Value = 'Monero'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Cells.Item(26, 4)
$c.Value = "'145.55"
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +2.17%  '

$ws.Cells.Item(27, 2).Value = 'Toncoin'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Cells.Item(27, 4)
$c.Value = "'1.654"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -0.05%  '

$ws.Cells.Item(28, 2).Value = 'EthereumClassic'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Cells.Item(28, 4)
$c.Value = "'17.19"
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +0.61%  '

$ws.Cells.Item(29, 2).Value = 'BitcoinCash'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c = $ws.Cells.Item(29, 4)
$c.Value = "'112.80"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +3.80%  '

$ws.Cells.Item(30, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Cells.Item(30, 4)
$c.Value = "'4.302"
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +0.31%  '

$ws.Cells.Item(31, 2).Value = 'Filecoin'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Cells.Item(31, 4)
$c.Value = "'4.270"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +1.02%  '

$ws.Cells.Item(32, 2).Value = 'Stellar'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Cells.Item(32, 4)
$c.Value = "'0.08834"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  +1.20%  '

$ws.Cells.Item(33, 2).Value = 'Hedera'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Cells.Item(33, 4)
$c.Value = "'0.04940"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  +2.98%  '

$ws.Cells.Item(34, 2).Value = 'ARBITRUM'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Cells.Item(34, 4)
$c.Value = "'1.157"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +1.64%  '

$ws.Cells.Item(35, 2).Value = 'ImmutableX'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Cells.Item(35, 4)
$c.Value = "'0.7204"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  -0.25%  '

$c = $ws.Cells.Item(36, 4)
$c.Value = "'2.881"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +1.54%  '

$ws.Cells.Item(37, 2).Value = 'MXToken'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Cells.Item(37, 4)
$c.Value = "'3.179"
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +1.46%  '

$ws.Cells.Item(38, 2).Value = 'VeChain'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Cells.Item(38, 4)
$c.Value = "'0.01834"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  -0.08%  '

$ws.Cells.Item(39, 2).Value = 'RenderToken'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Cells.Item(39, 4)
$c.Value = "'2.297"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -2.63%  '

$ws.Cells.Item(40, 2).Value = 'TheSandbox'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$c = $ws.Cells.Item(40, 4)
$c.Value = "'0.5042"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  -1.00%  '

$ws.Cells.Item(41, 2).Value = 'Quant'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c = $ws.Cells.Item(41, 4)
$c.Value = "'115.73"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -0.29%  '

$ws.Cells.Item(42, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Cells.Item(42, 4)
$c.Value = "'0.8952"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -4.90%  '

$ws.Cells.Item(43, 2).Value = 'FraxShare'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Cells.Item(43, 4)
$c.Value = "'6.119"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -1.01%  '

$ws.Cells.Item(44, 2).Value = 'Aptos'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c = $ws.Cells.Item(44, 4)
$c.Value = "'7.935"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +0.84%  '

$ws.Cells.Item(45, 2).Value = 'PaxDollar'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$c = $ws.Cells.Item(45, 4)
$c.Value = "'1.001"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +0.30%  '

$ws.Cells.Item(46, 2).Value = 'Decentraland'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$c = $ws.Cells.Item(46, 4)
$c.Value = "'0.4369"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  -1.32%  '

$ws.Cells.Item(47, 2).Value = 'Algorand'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Cells.Item(47, 4)
$c.Value = "'0.1318"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -2.72%  '

$ws.Cells.Item(48, 2).Value = 'EnergySwap'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Cells.Item(48, 4)
$c.Value = "'9.253"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -0.01%  '

$ws.Cells.Item(49, 2).Value = 'Elrond'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$c = $ws.Cells.Item(49, 4)
$c.Value = "'35.87"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -0.77%  '

$ws.Cells.Item(50, 2).Value = 'Cronos'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Cells.Item(50, 4)
$c.Value = "'0.05935"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +0.45%  '

$ws.Cells.Item(51, 2).Value = 'Aave'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Cells.Item(51, 4)
$c.Value = "'60.72"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +0.93%  '
